$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('AK_1.png', 0.458, 0.002, 0.795, 0.533, 'Akhlak Kamiswara', 'Benar')
    ,@('AK_2.png', 0.857, 0.003, 0.853, 0.867, 'Akhlak Kamiswara', 'Benar')
    ,@('AK_3.png', 0.841, 0.003, 0.856, 0.933, 'Akhlak Kamiswara', 'Benar')
    ,@('AK_4.png', 0.479, 0.002, 0.833, 0.533, 'Akhlak Kamiswara', 'Benar')
    ,@('AK_5.png', 0.442, 0.001, 0.854, 0.6, 'Akhlak Kamiswara', 'Benar')
    ,@('MIB_1.png', 1.543, 0.005, 0.661, 0.867, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_2.png', 1.162, 0.004, 0.742, 0.867, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_3.png', 1.741, 0.006, 0.952, 1, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_4.png', 1.162, 0.004, 0.746, 0.867, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('MIB_5.png', 1.492, 0.005, 0.714, 0.867, 'Muhammad Iqbal Baqi', 'Benar')
    ,@('AAH_1.png', 0.846, 0.003, 0.773, 0.933, 'Andrea Ayunove Hutami', 'Benar')
    ,@('AAH_2.png', 1.108, 0.004, 0.932, 1, 'Andrea Ayunove Hutami', 'Benar')
    ,@('AAH_3.png', 0.901, 0.003, 0.787, 1, 'Andrea Ayunove Hutami', 'Benar')
    ,@('TI_1.png', 0.96, 0.003, 0.792, 0.6, 'Toni Ismail', 'Benar')
    ,@('TI_2.png', 0.944, 0.003, 0.831, 0.733, 'Toni Ismail', 'Benar')
    ,@('TI_3.png', 0.75, 0.002, 0.934, 0.6, 'Toni Ismail', 'Benar')
    ,@('TI_4.png', 0.698, 0.002, 0.822, 0.8, 'Toni Ismail', 'Benar')
    ,@('TI_5.png', 1.065, 0.004, 0.87, 1, 'Toni Ismail', 'Benar')
    ,@('RAS_1.png', 0.637, 0.002, 0.844, 0.4, 'Tidak Diketahui', 'Salah')
    ,@('RAS_2.png', 1.032, 0.003, 0.74, 0.933, 'Ridha Ayu Salsabila', 'Benar')
    ,@('RAS_3.png', 0.543, 0.002, 0.781, 0.467, 'Tidak Diketahui', 'Salah')
    ,@('RAS_4.png', 1.25, 0.004, 0.693, 0.4, 'Tidak Diketahui', 'Salah')
    ,@('RAS_5.png', 1.177, 0.004, 0.78, 1, 'Ridha Ayu Salsabila', 'Benar')
    ,@('RR_1.png', 1.424, 0.005, 0.87, 0.867, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_2.png', 1.404, 0.005, 0.888, 0.867, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_3.png', 1.094, 0.004, 0.517, 0.733, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_4.png', 1.361, 0.005, 0.875, 0.933, 'Rafiqo Rapitasari', 'Benar')
    ,@('RR_5.png', 1.474, 0.005, 0.877, 0.933, 'Rafiqo Rapitasari', 'Benar')
    ,@('AR_1.png', 0.866, 0.003, 0.799, 1, 'Arizli Romadhon', 'Benar')
    ,@('GA_1.png', 1.548, 0.005, 0.879, 1, 'Gege Ardiyansyah', 'Benar')
    ,@('GA_2.png', 0.658, 0.002, 0.783, 0.333, 'Tidak Diketahui', 'Salah')
    ,@('GA_3.png', 0.734, 0.002, 0.761, 0.267, 'Tidak Diketahui', 'Salah')
    ,@('FY_1.png', 1.122, 0.004, 0.834, 0.333, 'Tidak Diketahui', 'Salah')
    ,@('FY_2.png', 1.628, 0.005, 0.861, 0.467, 'Tidak Diketahui', 'Salah')
    ,@('FY_3.png', 1.448, 0.005, 0.873, 0.533, 'Fanny Yusuf', 'Benar')
    ,@('FY_4.png', 1.256, 0.004, 0.802, 0.467, 'Tidak Diketahui', 'Salah')
    ,@('TO_1.png', 0.794, 0.003, 0.792, 0.8, 'Tiara Oktavian', 'Benar')
    ,@('TO_2.png', 1.058, 0.004, 0.817, 1, 'Tiara Oktavian', 'Benar')
    ,@('TO_3.png', 0.808, 0.003, 0.815, 1, 'Tiara Oktavian', 'Benar')
    ,@('TO_4.png', 3.35, 0.011, 0.514, 0.467, 'Tidak Diketahui', 'Salah')
    ,@('TO_5.png', 3.1, 0.01, 0.469, 0.467, 'Tidak Diketahui', 'Salah')
    ,@('TD_1.png', 2.326, 0.008, 0.42, 0.4, 'Tidak Diketahui', 'Benar')
    ,@('TD_2.png', 2.461, 0.008, 0.434, 0.4, 'Tidak Diketahui', 'Benar')
    ,@('TD_3.png', 1.038, 0.003, 0.716, 0.267, 'Tidak Diketahui', 'Benar')
    ,@('TD_4.png', 0.992, 0.003, 0.7, 0.267, 'Tidak Diketahui', 'Benar')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
